# Update the 25 "three-digit number x one-digit number" practice
# problems in the worksheet table. Each cell holds exactly one run of
# text like "742×3=" which we replace with the new expression.
#
# Note: one pair introduces a value collision - "692×7=" is replaced
# with "844×3=", while an earlier cell that already equals "844×3="
# is itself replaced with "659×6=" before that happens. Reusing a
# single Range (and letting Find.Execute advance it forward through
# the document on each call, in document order) guarantees each
# substitution only touches the intended occurrence and never
# re-matches a value that was just written by a later replacement.

$d = $word.ActiveDocument
$rng = $d.Content

$replacements = @(
    @("742×3=", "237×7="),
    @("246×4=", "621×5="),
    @("844×3=", "659×6="),
    @("724×8=", "523×4="),
    @("147×3=", "206×9="),
    @("567×9=", "795×4="),
    @("879×3=", "773×5="),
    @("999×8=", "969×2="),
    @("742×5=", "198×5="),
    @("571×3=", "646×5="),
    @("312×5=", "219×3="),
    @("457×3=", "106×7="),
    @("651×9=", "548×8="),
    @("301×4=", "826×4="),
    @("307×6=", "526×7="),
    @("658×4=", "832×8="),
    @("692×7=", "844×3="),
    @("451×7=", "607×6="),
    @("619×5=", "639×4="),
    @("187×6=", "880×6="),
    @("418×9=", "865×3="),
    @("237×4=", "597×7="),
    @("940×6=", "183×6="),
    @("968×6=", "494×9="),
    @("140×7=", "455×6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
